$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")
$lo = $ws.ListObjects.Item("Snippets")

# Add two new rows to the Snippets table for the new "Get word count" snippet
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# Fill in the new data, following the same column order used by the original author
$ws.Range("C34").Value = "word-paragraph-get-word-count"
$ws.Range("C35").Value = "word-paragraph-get-word-count"
$ws.Range("B35").Value = "paragraphs"
$ws.Range("D34").Value = "run"
$ws.Range("D35").Value = "run"
$ws.Range("B34").Value = "search"
$ws.Range("A34").Value = "Body"
$ws.Range("A35").Value = "Body"

# Reflect the final selection left by the editor
$ws.Range("D34").Select()
